$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Column C (rows 2-45) holds the "Förändrad" (changed/modified) date, stored
# as an Excel serial date number. Bump each from 45205 (2023-10-06) to
# 45206 (2023-10-07), matching the commit's automatic date update.
for ($row = 2; $row -le 45; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45205) {
        $cell.Value2 = 45206
    }
}
